$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H2").Value = 161.18182
$ws.Range("I2").Value = 66.333336
$ws.Range("J2").Value = 196.75
$ws.Range("K2").Value = 66.333336
$ws.Range("L2").Value = 196.75
$ws.Range("M2").Value = 46.666664
$ws.Range("N2").Value = -422.75

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H21").Value = 14998.5
$ws.Range("J21").Value = 14998.5
$ws.Range("L21").Value = 14998.5
$ws.Range("N21").Value = -15934.5

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H23").Value = 14998.5
$ws.Range("J23").Value = 14998.5
$ws.Range("L23").Value = 14998.5
$ws.Range("N23").Value = -15466.5

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H32").Value = 1887.8889
$ws.Range("I32").Value = 1999.8334
$ws.Range("K32").Value = 1999.8334
$ws.Range("M32").Value = -1673.8334

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H39").Value = 177.33333
$ws.Range("I39").Value = 93.888885
$ws.Range("J39").Value = 302.5
$ws.Range("K39").Value = 281.666655
$ws.Range("L39").Value = 907.5
$ws.Range("M39").Value = 14.33334500000001
$ws.Range("N39").Value = -1499.5

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H57").Value = 75660
$ws.Range("J57").Value = 75660
$ws.Range("L57").Value = 226980
$ws.Range("N57").Value = -227978

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H98").Value = 827.2222
$ws.Range("I98").Value = 827.2222
$ws.Range("K98").Value = 827.2222
$ws.Range("M98").Value = 670.7778

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H122").Value = 827.2222
$ws.Range("I122").Value = 827.2222
$ws.Range("K122").Value = 2481.6666
$ws.Range("M122").Value = -31.66660000000002

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H138").Value = 3384.2
$ws.Range("I138").Value = 2661.35
$ws.Range("J138").Value = 4829.9
$ws.Range("K138").Value = 7984.049999999999
$ws.Range("L138").Value = 14489.7
$ws.Range("M138").Value = -2844.049999999999
$ws.Range("N138").Value = -24769.7

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H92").Value = 65000
$ws.Range("J92").Value = 65000
$ws.Range("L92").Value = 65000
$ws.Range("N92").Value = -69992

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H134").Value = 77999.336
$ws.Range("J134").Value = 77999.336
$ws.Range("L134").Value = 77999.336
$ws.Range("N134").Value = -88139.336

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H22").Value = 2103.375
$ws.Range("I22").Value = 269.5
$ws.Range("K22").Value = 269.5
$ws.Range("M22").Value = -96.5

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H59").Value = 96350
$ws.Range("J59").Value = 96350
$ws.Range("L59").Value = 96350
$ws.Range("N59").Value = -98044

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H60").Value = 105000
$ws.Range("J60").Value = 105000
$ws.Range("L60").Value = 105000
$ws.Range("N60").Value = -106198

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H106").Value = 70000
$ws.Range("J106").Value = 70000
$ws.Range("L106").Value = 70000
$ws.Range("N106").Value = -72524

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H132").Value = 98994.28999999999
$ws.Range("J132").Value = 98994.28999999999
$ws.Range("L132").Value = 98994.28999999999
$ws.Range("N132").Value = -109114.29

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 5558431
$ws.Range("I134").Value = 2974.5173
$ws.Range("J134").Value = 166666670
$ws.Range("K134").Value = 8923.5519
$ws.Range("L134").Value = 500000010
$ws.Range("M134").Value = -6388.5519
$ws.Range("N134").Value = -500005080

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H28").Value = 72949.836
$ws.Range("J28").Value = 27539.8
$ws.Range("L28").Value = 27539.8
$ws.Range("N28").Value = -28029.8

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 3761084.5
$ws.Range("I31").Value = 4050206.5
$ws.Range("J31").Value = 2499
$ws.Range("K31").Value = 4050206.5
$ws.Range("L31").Value = 2499
$ws.Range("M31").Value = -4049911.5
$ws.Range("N31").Value = -3089

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H34").Value = 3761084.5
$ws.Range("I34").Value = 4050206.5
$ws.Range("J34").Value = 2499
$ws.Range("K34").Value = 4050206.5
$ws.Range("L34").Value = 2499
$ws.Range("M34").Value = -4050004.5
$ws.Range("N34").Value = -2903

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H52").Value = 62356
$ws.Range("J52").Value = 65445
$ws.Range("L52").Value = 65445
$ws.Range("N52").Value = -66033

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 16111961
$ws.Range("I58").Value = 27782350
$ws.Range("J58").Value = 7359169.5
$ws.Range("K58").Value = 27782350
$ws.Range("L58").Value = 7359169.5
$ws.Range("M58").Value = -27782147
$ws.Range("N58").Value = -7359575.5

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H75").Value = 28333.334
$ws.Range("J75").Value = 20000
$ws.Range("L75").Value = 20000
$ws.Range("N75").Value = -21996

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H78").Value = 28333.334
$ws.Range("J78").Value = 20000
$ws.Range("L78").Value = 60000
$ws.Range("N78").Value = -69984

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H88").Value = 12508
$ws.Range("J88").Value = 12508
$ws.Range("L88").Value = 12508
$ws.Range("N88").Value = -13320

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H91").Value = 12508
$ws.Range("J91").Value = 12508
$ws.Range("L91").Value = 12508
$ws.Range("N91").Value = -15316

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H136").Value = 16111961
$ws.Range("I136").Value = 27782350
$ws.Range("J136").Value = 7359169.5
$ws.Range("K136").Value = 83347050
$ws.Range("L136").Value = 22077508.5
$ws.Range("M136").Value = -83344500
$ws.Range("N136").Value = -22082608.5

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 147.76471
$ws.Range("J2").Value = 143.78572
$ws.Range("L2").Value = 862.71432
$ws.Range("N2").Value = -1088.71432

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 8173787
$ws.Range("I4").Value = 15800631
$ws.Range("K4").Value = 47401893
$ws.Range("M4").Value = -47401781

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H23").Value = 171.33333
$ws.Range("J23").Value = 227.88889
$ws.Range("L23").Value = 683.6666700000001
$ws.Range("N23").Value = -1153.66667

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H105").Value = 0
$ws.Range("J105").Value = 0
$ws.Range("L105").Value = 0
$ws.Range("N105").ClearContents()

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H119").Value = 0
$ws.Range("J119").Value = 0
$ws.Range("L119").Value = 0
$ws.Range("N119").ClearContents()

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H136").Value = 55108.332
$ws.Range("J136").Value = 55108.332
$ws.Range("L136").Value = 165324.996
$ws.Range("N136").Value = -170424.996

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 3693.5454
$ws.Range("I46").Value = 1723.0769
$ws.Range("K46").Value = 1723.0769
$ws.Range("M46").Value = -1535.0769

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H93").Value = 11913.333
$ws.Range("I93").Value = 9690.143
$ws.Range("K93").Value = 9690.143
$ws.Range("M93").Value = -8442.143

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H133").Value = 82835.8
$ws.Range("J133").Value = 82835.8
$ws.Range("L133").Value = 82835.8
$ws.Range("N133").Value = -87895.8

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H140").Value = 99998
$ws.Range("J140").Value = 99998
$ws.Range("L140").Value = 99998
$ws.Range("N140").Value = -110358

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H76").Value = 31000
$ws.Range("J76").Value = 40000
$ws.Range("L76").Value = 40000
$ws.Range("N76").Value = -40630

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H79").Value = 31000
$ws.Range("J79").Value = 40000
$ws.Range("L79").Value = 40000
$ws.Range("N79").Value = -42184

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 8336651
$ws.Range("I132").Value = 11114107
$ws.Range("K132").Value = 33342321
$ws.Range("M132").Value = -33339791

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H137").Value = 61738.332
$ws.Range("J137").Value = 61738.332
$ws.Range("L137").Value = 61738.332
$ws.Range("N137").Value = -71938.332
